$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that ends with:
# "...Some of the analyst's other forecasts are summarized in the
#  table below:" -- a new "pander(table_forecasts)" source-code
# paragraph must be inserted immediately after it (and before the
# table that follows).
# ------------------------------------------------------------------
$anchorRange = $d.Content
$found = $anchorRange.Find.Execute("Some of the analyst", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the anchor paragraph (forecast summary sentence)."
}
$anchorRange.Expand(4) | Out-Null   # wdParagraph -> grow the found hit to the whole paragraph

$anchorStart = $anchorRange.Start
$anchorEnd = $anchorRange.End

# Resolve the Paragraphs-collection index for this exact paragraph
# (collection indices inside table regions are not reliably ordered,
# so match on both Start and End to be safe).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pr = $d.Paragraphs.Item($i).Range
    if ($pr.Start -eq $anchorStart -and $pr.End -eq $anchorEnd) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not resolve the anchor paragraph's index."
}

# ------------------------------------------------------------------
# Insert a new, empty paragraph right after the anchor paragraph.
# ------------------------------------------------------------------
$paraEnd = $d.Paragraphs.Item($anchorIndex).Range.End
$insertionPoint = $d.Range($paraEnd - 1, $paraEnd - 1)
$insertionPoint.InsertParagraphAfter()

$newParaIndex = $anchorIndex + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newPara.Style = "Source Code"

# ------------------------------------------------------------------
# First run: "pander" styled as FunctionTok
# ------------------------------------------------------------------
$r = $newPara.Range
$r.Text = "pander"

$runStart = $r.Start
$runEnd = $r.End
$firstRunRange = $d.Range($runStart, $runEnd - 1)
$firstRunRange.Style = "FunctionTok"

# ------------------------------------------------------------------
# Second run: "(table_forecasts)" styled as NormalTok
# ------------------------------------------------------------------
$insertAfterFirst = $d.Range($runEnd - 1, $runEnd - 1)
$insertAfterFirst.InsertAfter("(table_forecasts)")

$paraRangeEnd = $newPara.Range.End
$secondRunRange = $d.Range($runEnd - 1, $paraRangeEnd - 1)
$secondRunRange.Style = "NormalTok"

Write-Output "Inserted source-code paragraph after paragraph #${anchorIndex}: [$($newPara.Range.Text)]"
